{"js": "// Apply the OCR-style text corrections and remove the two inline-picture\n// paragraphs, matching the target OOXML diff.\n\nconst body = context.document.body;\n\n// --- 1) Targeted text fixes (exact substring search + literal replace) ---\nconst replacements = [\n  [\"Section A (28 x 2 marks = 56 marks)\", \"Section A (28 x2 marks = 56 marks)\"],\n  [\n    \"For each. question from 1 to 28, four options are given. One of them is the correct answer.\",\n    \"For each question from 1 to 28, four options are given. One of them is the correct answer.\",\n  ],\n  [\n    \"Make your choice (1, 2, 3.0r4) and shade-your answer on the Optical Answer Sheet provided.\",\n    \"Make your choice (1, 2, 3 or4) and shade-your answer on the Optical Answer Sheet provided.\",\n  ],\n  [\n    \"1. The two graphs below show how the light intensity and the depth of the pond affect the\",\n    \"1. The two grapns below show now the lignt intensity and the depth of ine pond affect the\",\n  ],\n  [\n    \"(1) The tight intensity aoes not affect ihe rate of photosynthesis.\",\n    \"(7) dhe ugnt intensity aoes not affect ine rate of photosynthesis.\",\n  ],\n  [\n    \"(3) The rate of photosynthesis does not depend on the deptfof the pond.\",\n    \"(3)  Therate of photosynthesis does not depend on the depth of the pond.\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n// --- 2) Remove the two paragraphs that each hold a single inline picture ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst toDelete = [];\nfor (const para of paragraphs.items) {\n  const pics = para.inlinePictures;\n  pics.load(\"items\");\n  await context.sync();\n  if (pics.items.length > 0) {\n    toDelete.push(para);\n  }\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\nawait context.sync();\n", "ps1": "# Apply the OCR-style text corrections and remove the two inline-picture\n# paragraphs, matching the target OOXML diff.\n\n$d = $word.ActiveDocument\n\n# --- 1) Targeted text fixes (exact literal find & replace) ---\n$pairs = @(\n    @(\"Section A (28 x 2 marks = 56 marks)\", \"Section A (28 x2 marks = 56 marks)\"),\n    @(\"For each. question from 1 to 28, four options are given. One of them is the correct answer.\", \"For each question from 1 to 28, four options are given. One of them is the correct answer.\"),\n    @(\"Make your choice (1, 2, 3.0r4) and shade-your answer on the Optical Answer Sheet provided.\", \"Make your choice (1, 2, 3 or4) and shade-your answer on the Optical Answer Sheet provided.\"),\n    @(\"1. The two graphs below show how the light intensity and the depth of the pond affect the\", \"1. The two grapns below show now the lignt intensity and the depth of ine pond affect the\"),\n    @(\"(1) The tight intensity aoes not affect ihe rate of photosynthesis.\", \"(7) dhe ugnt intensity aoes not affect ine rate of photosynthesis.\"),\n    @(\"(3) The rate of photosynthesis does not depend on the deptfof the pond.\", \"(3)  Therate of photosynthesis does not depend on the depth of the pond.\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# --- 2) Remove the two paragraphs that each hold a single inline picture ---\nfor ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {\n    $shape = $d.InlineShapes.Item($i)\n    $shape.Range.Paragraphs.Item(1).Range.Delete()\n}\n"}
